$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Adafruit IO reading appended as row 90 (duplicate of the last
# polled "temperature" data point, matching the feed's last entry).
$ws.Range("A90").Value = "2024-09-25T18:06:40Z"
$ws.Range("B90").Value = "temperature"

# The "Value" column stores numeric-looking readings as text (as every
# other row in this column already does). Force a text number format
# just long enough to write the value as a string, then clear the
# formatting again so the cell keeps the default style used elsewhere
# in the sheet.
$ws.Range("C90").NumberFormat = "@"
$ws.Range("C90").Value = "25"
$ws.Range("C90").ClearFormats()

$ws.Range("D90").Value = "N/A"
$ws.Range("E90").Value = "N/A"
$ws.Range("F90").Value = "N/A"
